$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 from 22:15 to 22:45
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 22:45"

# Swap the labels for the duplicate "Illes Balears" / "Illes Balears*" rows
$ws.Range("A26").Value = "Illes Balears"
$ws.Range("A27").Value = "Illes Balears*"

# Swap the labels for the duplicate "Huelva" / "Melilla" rows
$ws.Range("A52").Value = "Melilla"
$ws.Range("A53").Value = "Huelva"

# Swap the labels for the duplicate "Ceuta" / "La Palma" rows
$ws.Range("A57").Value = "La Palma"
$ws.Range("A58").Value = "Ceuta"
